# reference_table_sector.xlsx — "new esto file name"
#
# The only in-sheet / model-level change in the target revision is the
# worksheet being renamed from the default CJK placeholder "工作表1" to the
# descriptive name "sectors" (everything else in the upstream diff is
# session/host metadata — absPath, revisionPtr/documentId GUIDs, the
# workbookView x/y window coordinates, customXml part shuffling — that
# Excel itself stamps on save and isn't reachable through the workbook
# object model).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the single worksheet to match the new file naming convention.
$ws.Name = "sectors"

# Best-effort: nudge the window scroll position / placement to line up with
# the new view (topLeftCell A149 in the target). Harmless no-op on hosts
# that don't expose these, but mirrors the intent of the source edit.
$win = $excel.ActiveWindow
[void]$ws.Range("A149").Select()
$win.ScrollRow = 149
$win.ScrollColumn = 1
$win.Left = 28680
$win.Top = -120

# Restore the original active cell selection (the diff keeps C118 selected).
[void]$ws.Range("C118").Select()
